$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4000
$ws.Range("I17").Value = 4000
$ws.Range("K17").Value = 12000
$ws.Range("M17").Value = -11832
$ws.Range("H28").Value = 1857.3334
$ws.Range("I28").Value = 836
$ws.Range("J28").Value = 3900
$ws.Range("K28").Value = 836
$ws.Range("L28").Value = 3900
$ws.Range("M28").Value = -351
$ws.Range("N28").Value = -4870
$ws.Range("H33").Value = 271.63635
$ws.Range("I33").Value = 206.3158
$ws.Range("K33").Value = 206.3158
$ws.Range("M33").Value = 22.6842
$ws.Range("H58").Value = 104.75
$ws.Range("I58").Value = 104.75
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 314.25
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -164.25
$ws.Range("H62").Value = 3674.75
$ws.Range("I62").Value = 3233
$ws.Range("K62").Value = 3233
$ws.Range("M62").Value = -2609
$ws.Range("H64").Value = 10003.75
$ws.Range("I64").Value = 10007.5
$ws.Range("K64").Value = 10007.5
$ws.Range("M64").Value = -9759.5
$ws.Range("H65").Value = 3674.75
$ws.Range("I65").Value = 3233
$ws.Range("K65").Value = 16165
$ws.Range("M65").Value = -13045
$ws.Range("H67").Value = 10003.75
$ws.Range("I67").Value = 10007.5
$ws.Range("K67").Value = 10007.5
$ws.Range("M67").Value = -9149.5
$ws.Range("H107").Value = 52988.555
$ws.Range("I107").Value = 67594.64
$ws.Range("J107").Value = 1867.25
$ws.Range("K107").Value = 67594.64
$ws.Range("L107").Value = 1867.25
$ws.Range("M107").Value = -65674.64
$ws.Range("N107").Value = -5707.25
$ws.Range("H132").Value = 2146.3333
$ws.Range("I132").Value = 2514.4285
$ws.Range("J132").Value = 1410.1428
$ws.Range("K132").Value = 7543.2855
$ws.Range("L132").Value = 4230.428400000001
$ws.Range("M132").Value = -5013.2855
$ws.Range("N132").Value = -9290.428400000001
$ws.Range("H137").Value = 2111.75
$ws.Range("I137").Value = 1628.4
$ws.Range("J137").Value = 2917.3333
$ws.Range("K137").Value = 4885.200000000001
$ws.Range("L137").Value = 8751.999899999999
$ws.Range("M137").Value = -2335.200000000001
$ws.Range("N137").Value = -13851.9999
$ws.Range("H138").Value = 7075.263
$ws.Range("I138").Value = 7049.75
$ws.Range("J138").Value = 7082.067
$ws.Range("K138").Value = 21149.25
$ws.Range("L138").Value = 21246.201
$ws.Range("M138").Value = -16009.25
$ws.Range("N138").Value = -31526.201
$ws.Range("H141").Value = 916.6667
$ws.Range("I141").Value = 916.6667
$ws.Range("K141").Value = 2750.0001
$ws.Range("M141").Value = 2429.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 20750.5
$ws.Range("I12").Value = 20251.5
$ws.Range("J12").Value = 21000
$ws.Range("K12").Value = 20251.5
$ws.Range("L12").Value = 21000
$ws.Range("M12").Value = -20078.5
$ws.Range("N12").Value = -21346
$ws.Range("H45").Value = 1855.5
$ws.Range("I45").Value = 1855.5
$ws.Range("K45").Value = 1855.5
$ws.Range("M45").Value = -1478.5
$ws.Range("H61").Value = 933
$ws.Range("I61").Value = 933
$ws.Range("K61").Value = 933
$ws.Range("M61").Value = -721
$ws.Range("H102").Value = 35002784
$ws.Range("I102").Value = 2003342
$ws.Range("K102").Value = 2003342
$ws.Range("M102").Value = -2001720
$ws.Range("H122").Value = 8011.8423
$ws.Range("J122").Value = 9601.923000000001
$ws.Range("L122").Value = 28805.769
$ws.Range("N122").Value = -33705.769
$ws.Range("H132").Value = 3856.5715
$ws.Range("I132").Value = 2999.2
$ws.Range("K132").Value = 8997.599999999999
$ws.Range("M132").Value = -6467.599999999999
$ws.Range("H136").Value = 933
$ws.Range("I136").Value = 933
$ws.Range("K136").Value = 2799
$ws.Range("M136").Value = -249
$ws.Range("H105").Value = 1481.3334
$ws.Range("I105").Value = 1000
$ws.Range("J105").Value = 2444
$ws.Range("K105").Value = 1000
$ws.Range("L105").Value = 2444
$ws.Range("M105").Value = 747
$ws.Range("N105").Value = -5938

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3122.1428
$ws.Range("I134").Value = 3143
$ws.Range("K134").Value = 9429
$ws.Range("M134").Value = -6894

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 227.09091
$ws.Range("I7").Value = 104.6
$ws.Range("K7").Value = 104.6
$ws.Range("M7").Value = 8.400000000000006
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H41").Value = 33026.25
$ws.Range("J41").Value = 36978.715
$ws.Range("L41").Value = 36978.715
$ws.Range("N41").Value = -37834.715
$ws.Range("H50").Value = 45000
$ws.Range("J50").Value = 45000
$ws.Range("L50").Value = 45000
$ws.Range("N50").Value = -46250
$ws.Range("H99").Value = 1002782.4
$ws.Range("I99").Value = 1670000
$ws.Range("J99").Value = 1956
$ws.Range("K99").Value = 1670000
$ws.Range("L99").Value = 1956
$ws.Range("M99").Value = -1668502
$ws.Range("N99").Value = -4952
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("N122").Value = 0
$ws.Range("H126").Value = 1002782.4
$ws.Range("I126").Value = 1670000
$ws.Range("J126").Value = 1956
$ws.Range("K126").Value = 5010000
$ws.Range("L126").Value = 5868
$ws.Range("M126").Value = -5007530
$ws.Range("N126").Value = -10808
$ws.Range("H132").Value = 199.66667
$ws.Range("I132").Value = 199.66667
$ws.Range("K132").Value = 599.00001
$ws.Range("M132").Value = 1930.99999
$ws.Range("H134").Value = 1637.5
$ws.Range("I134").Value = 1565
$ws.Range("K134").Value = 4695
$ws.Range("M134").Value = -2160
$ws.Range("H141").Value = 959582.5
$ws.Range("J141").Value = 959582.5
$ws.Range("L141").Value = 959582.5
$ws.Range("N141").Value = -969942.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 873.4167
$ws.Range("I6").Value = 47.625
$ws.Range("J6").Value = 2525
$ws.Range("K6").Value = 142.875
$ws.Range("L6").Value = 7575
$ws.Range("M6").Value = -29.875
$ws.Range("N6").Value = -7801
$ws.Range("H113").Value = 999
$ws.Range("J113").Value = 999
$ws.Range("L113").Value = 2997
$ws.Range("N113").Value = -7337
$ws.Range("H121").Value = 29491.5
$ws.Range("I121").Value = 449.5
$ws.Range("J121").Value = 35299.9
$ws.Range("K121").Value = 1348.5
$ws.Range("L121").Value = 105899.7
$ws.Range("M121").Value = -38.5
$ws.Range("N121").Value = -108519.7
$ws.Range("H129").Value = 2744.0715
$ws.Range("I129").Value = 2988.4285
$ws.Range("J129").Value = 2499.7144
$ws.Range("K129").Value = 8965.2855
$ws.Range("L129").Value = 7499.1432
$ws.Range("M129").Value = -3965.2855
$ws.Range("N129").Value = -17499.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 116.65
$ws.Range("I2").Value = 95.53333000000001
$ws.Range("K2").Value = 95.53333000000001
$ws.Range("M2").Value = 17.46666999999999
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").ClearContents()
$ws.Range("N114").Value = 0
$ws.Range("H122").Value = 3146.2
$ws.Range("I122").Value = 3289.6667
$ws.Range("K122").Value = 9869.000100000001
$ws.Range("M122").Value = -7419.000100000001
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").ClearContents()
$ws.Range("N128").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I14").Value = 0
$ws.Range("J14").ClearContents()
$ws.Range("K14").Value = 0
$ws.Range("L14").ClearContents()
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -1344
$ws.Range("H16").Value = 997.55554
$ws.Range("I16").Value = 1013.1667
$ws.Range("K16").Value = 1013.1667
$ws.Range("M16").Value = -843.1667
$ws.Range("H22").Value = 406.25
$ws.Range("I22").Value = 312.5
$ws.Range("K22").Value = 312.5
$ws.Range("M22").Value = -17.5
$ws.Range("H27").Value = 406.25
$ws.Range("I27").Value = 312.5
$ws.Range("K27").Value = 312.5
$ws.Range("M27").Value = -205.5
$ws.Range("H35").Value = 12696.417
$ws.Range("I35").Value = 7060.3335
$ws.Range("J35").Value = 18332.5
$ws.Range("K35").Value = 7060.3335
$ws.Range("L35").Value = 18332.5
$ws.Range("M35").Value = -6724.3335
$ws.Range("N35").Value = -19004.5
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = 0
$ws.Range("H134").Value = 24999
$ws.Range("J134").Value = 24999
$ws.Range("L134").Value = 24999
$ws.Range("N134").Value = -35139
$ws.Range("H136").Value = 4999
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").ClearContents()
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = 0
$ws.Range("H132").Value = 4499
$ws.Range("I132").Value = 1000
$ws.Range("K132").Value = 3000
$ws.Range("M132").Value = -470
